# Update "想去人数" (attendance interest count) values across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1130
$ws.Range("F4").Value = 1205
$ws.Range("F6").Value = 161
$ws.Range("F8").Value = 292
$ws.Range("F11").Value = 28098
$ws.Range("F12").Value = 3270
$ws.Range("F14").Value = 235
$ws.Range("F15").Value = 452
$ws.Range("F16").Value = 8
$ws.Range("F17").Value = 40
$ws.Range("F19").Value = 305
$ws.Range("F20").Value = 586
$ws.Range("F21").Value = 262
$ws.Range("F22").Value = 243
$ws.Range("F25").Value = 26
$ws.Range("F27").Value = 190
$ws.Range("F28").Value = 87
$ws.Range("F29").Value = 500
$ws.Range("F30").Value = 67
$ws.Range("F32").Value = 590
$ws.Range("F33").Value = 234
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 364
$ws.Range("F7").Value = 762
$ws.Range("F22").Value = 4228
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 243
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 243
$ws.Range("F9").Value = 364
$ws.Range("F11").Value = 762
$ws.Range("F12").Value = 1130
$ws.Range("F13").Value = 1205
$ws.Range("F14").Value = 161
$ws.Range("F16").Value = 292
$ws.Range("F27").Value = 235
$ws.Range("F30").Value = 452
$ws.Range("F31").Value = 8
$ws.Range("F34").Value = 305
$ws.Range("F35").Value = 586
$ws.Range("F36").Value = 262
$ws.Range("F39").Value = 26
$ws.Range("F42").Value = 190
$ws.Range("F43").Value = 87
$ws.Range("F46").Value = 67
$ws.Range("F48").Value = 590
$ws.Range("F49").Value = 234
